# Apply "Added caching data feature for exchanges" changes.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("ListOfValues")

# On Sheet1, the Exchange column (B) values change from "ByBit" to "Binance"
$ws1.Range("B2").Value = "Binance"
$ws1.Range("B3").Value = "Binance"

# On ListOfValues, the Strategy list entry "Early MACD" becomes "EarlyMACD"
$ws2.Range("B3").Value = "EarlyMACD"

# Update selections to match the saved state in the diff
$ws1.Activate()
$ws1.Range("L3").Select()

$ws2.Activate()
$ws2.Range("E9").Select()

# Re-activate Sheet1 as the active sheet (tabSelected="1" on sheet1)
$ws1.Activate()
